$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns B,C,D,F,G,L,M,O
$data = @(
    @(2, 1.095848521699395, 0.3766966690989761, 0.02983574819897683, 0.7669624737050498, 0.00243212134707261, 0.2701266419800135, 0.2543820577602602, 2.649840988883795),
    @(3, 0.9896832470606114, 0.3637491913727899, 0.02832589650306261, 0.7667236685179191, 0.002435121979148858, 0.2689680485931163, 0.2385282370809136, 2.664703294676002),
    @(4, 0.9245442686059278, 0.3557893218381025, 0.02739473749461752, 0.7671982436097551, 0.00243706301914252, 0.2684068073203605, 0.228864422130755, 2.675951353872847),
    @(5, 0.8980128366425788, 0.3525433650557659, 0.02701427779389576, 0.7675475738428901, 0.002437878887630242, 0.2682158706442763, 0.2249442616060122, 2.681067784640021),
    @(6, 0.8936081540948067, 0.3520042492646098, 0.02695104286446082, 0.7676149901686671, 0.002438015866765183, 0.2681864477861211, 0.2242944102309323, 2.681949512555917),
    @(7, 0.9241864009583196, 0.3557455544399204, 0.02738961051468891, 0.7672023237760897, 0.0024370739214153, 0.2684040793105424, 0.2288114806458239, 2.676018200044609),
    @(8, 1.059233815033338, 0.3722346343712388, 0.02931601884390034, 0.7667510387421288, 0.002433135537625974, 0.2696960081912678, 0.2489011408707356, 2.65452438763225),
    @(9, 1.324384722473781, 0.4044789424853263, 0.03306013259597762, 0.7708084382374523, 0.002426191611163806, 0.2734204153761297, 0.2888498465199234, 2.62926288531412),
    @(10, 1.51934265579348, 0.4281014261370046, 0.03578939659328029, 0.7768227769668101, 0.002421560139978345, 0.2768830601651615, 0.3185318990945305, 2.621066124280219),
    @(11, 1.608058796412593, 0.4388309081520276, 0.03702612717425069, 0.7802219931706276, 0.00241955424637246, 0.2786161388830664, 0.3321061649565067, 2.619602122042181),
    @(12, 1.641656320943923, 0.4428912662944526, 0.03749372829816622, 0.781604887578851, 0.002418809110888361, 0.2792951120056131, 0.3372565597314221, 2.619374538488273),
    @(13, 1.634420398769976, 0.4420169185815439, 0.03739305457945363, 0.781302796016476, 0.002418968947419706, 0.2791478738809587, 0.3361468834102368, 2.619409000039468),
    @(14, 1.610822840417825, 0.4391650112981722, 0.03706461162776264, 0.7803338452957433, 0.002419492654213504, 0.2786715436644442, 0.3325296894716487, 2.619576842519024),
    @(15, 1.596368955274272, 0.437417780271403, 0.03686333599626579, 0.7797528051362406, 0.002419815320167234, 0.2783827326867083, 0.3303153650114581, 2.619722243785532),
    @(16, 1.51354529881911, 0.4273998722248962, 0.03570847383978304, 0.776614003168973, 0.002421693255983648, 0.2767729757513706, 0.3176462153293897, 2.621207469844677),
    @(17, 1.462741976490463, 0.4212497750645525, 0.03499874738275821, 0.7748585476610117, 0.00242287112481881, 0.2758258760607788, 0.3098923346023668, 2.622699472481969),
    @(18, 1.433524048166646, 0.4177108655115376, 0.03459007896947952, 0.7739112672472714, 0.002423558113000392, 0.2752959928107686, 0.3054392730586528, 2.623770733955325),
    @(19, 1.423631886170881, 0.4165123965929638, 0.03445163400472495, 0.7736012441218776, 0.002423792350678509, 0.2751191366434824, 0.3039327111406323, 2.624170010347115),
    @(20, 1.468149798068168, 0.4219046245282527, 0.03507434601374371, 0.7750389573323702, 0.00242274475486679, 0.2759251582215967, 0.3107170497080602, 2.622518582855918),
    @(21, 1.61775395529412, 0.4400027601921011, 0.03716110309839138, 0.7806158502975364, 0.002419338437068593, 0.2788108375912515, 0.3335918741049397, 2.619518664773977),
    @(22, 1.715543669697638, 0.4518152789830481, 0.03852069889153142, 0.7848185394536529, 0.002417196422199151, 0.2808290670911617, 0.3486007588146336, 2.619463422110101),
    @(23, 1.663350598659974, 0.4455122432351857, 0.03779545238889881, 0.782524338841398, 0.002418331973199823, 0.2797398013318855, 0.3405849149660725, 2.619318190358314),
    @(24, 1.465704954285627, 0.4216085768341031, 0.03504016985618819, 0.7749572011577897, 0.002422801856210933, 0.2758802272260823, 0.3103441812055578, 2.622599698126095),
    @(25, 1.252624532376501, 0.3957669810866378, 0.03205095645453326, 0.7691795773028645, 0.002427987204045439, 0.2722853365068545, 0.2779840279621553, 2.634281830457326),
)

$colLetters = @("B", "C", "D", "F", "G", "L", "M", "O")

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $val = $entry[$i + 1]
        $ws.Range("$col$row").Value = $val
    }
}

Write-Host "Updated $($data.Length) rows"
